# Stocks_Excel.py now works using trial lists
# - Renames the existing sheet "Sheet" to "Sheet6" (keeps it first / active).
# - Adds a brand-new sheet named "Sheet" right after it.
# - Populates both sheets with the ticker/price trial-list data.

$wb = $excel.ActiveWorkbook

# The workbook currently has a single sheet named "Sheet" - it becomes "Sheet6".
$ws1 = $wb.ActiveSheet
$ws1.Name = "Sheet6"

# Add the new "Sheet" worksheet right after Sheet6.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet"

# Match the original sheet's outline/page properties on the new sheet.
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = $ws1.PageSetup.LeftMargin
$ws2.PageSetup.RightMargin = $ws1.PageSetup.RightMargin
$ws2.PageSetup.TopMargin = $ws1.PageSetup.TopMargin
$ws2.PageSetup.BottomMargin = $ws1.PageSetup.BottomMargin
$ws2.PageSetup.HeaderMargin = $ws1.PageSetup.HeaderMargin
$ws2.PageSetup.FooterMargin = $ws1.PageSetup.FooterMargin

# --- Sheet6 (first sheet): trial list of tickers only ---
$ws1.Range("A1").NumberFormat = "@"
$ws1.Range("A1").Value = "01/04/21"
$ws1.Range("C1").Value = "Ticker"
$ws1.Range("D1").Value = "Price"
$ws1.Range("C2").Value = "stock1"
$ws1.Range("C3").Value = "stock2"
$ws1.Range("C4").Value = "stock3"
$ws1.Range("C5").Value = "stock4"
$ws1.Range("C6").Value = "stock5"

# --- Sheet (second sheet): tickers with matching prices ---
$ws2.Range("A1").NumberFormat = "@"
$ws2.Range("A1").Value = "01/04/21"
$ws2.Range("C1").Value = "Ticker"
$ws2.Range("D1").Value = "Price"
$ws2.Range("C2").Value = "stock1"
$ws2.Range("D2").Value = "price1"
$ws2.Range("C3").Value = "stock2"
$ws2.Range("D3").Value = "price2"
$ws2.Range("C4").Value = "stock3"
$ws2.Range("D4").Value = "price3"
$ws2.Range("C5").Value = "stock4"
$ws2.Range("D5").Value = "price4"
$ws2.Range("C6").Value = "stock5"
$ws2.Range("D6").Value = "price5"

# --- Selections as saved in the original workbook ---
[void]$ws2.Range("A1").Select()
[void]$ws1.Activate()
[void]$ws1.Range("C13").Select()
